$wb = $excel.ActiveWorkbook

# Sheet 5: "具有相當價值之財產" (property of considerable value) -
# fund, bonds, otherbonds, antique. The sheet already has a 2-row table
# (A1:E2, name/quantity/owner/total columns) for the one "珠寶材料"
# (jewelry material) antique entry; extend it out to A1:L2 with the
# property_category/category/date/legislator_name/legislator_id/
# source_file/index columns used on every other sheet (see "存款").
$ws5 = $wb.Worksheets.Item(5)

# --- header row (row 1) ---
$ws5.Cells.Item(1, 2).Value = "name"
$ws5.Cells.Item(1, 3).Value = "quantity"
$ws5.Cells.Item(1, 4).Value = "owner"
$ws5.Cells.Item(1, 5).Value = "total"
$ws5.Cells.Item(1, 6).Value = "property_category"
$ws5.Cells.Item(1, 7).Value = "category"
$ws5.Cells.Item(1, 8).Value = "date"
$ws5.Cells.Item(1, 9).Value = "legislator_name"
$ws5.Cells.Item(1, 10).Value = "legislator_id"
$ws5.Cells.Item(1, 11).Value = "source_file"
$ws5.Cells.Item(1, 12).Value = "index"

# --- data row (row 2) ---
$ws5.Cells.Item(2, 2).Value = "珠寶材料"
$ws5.Cells.Item(2, 3).Value = 40
$ws5.Cells.Item(2, 4).Value = "黃素香"
$ws5.Cells.Item(2, 5).Value = "6000000(製作珠寶飾品之原材料（估計價值））"
$ws5.Cells.Item(2, 6).Value = "otherbonds"
$ws5.Cells.Item(2, 7).Value = "normal"
# Leading apostrophe forces this to stay a text value ("2012-04-25")
# instead of Excel auto-converting it to a date serial number.
$ws5.Cells.Item(2, 8).Value = "'2012-04-25"
$ws5.Cells.Item(2, 9).Value = "李桐豪"
$ws5.Cells.Item(2, 10).Value = 896
$ws5.Cells.Item(2, 11).Value = "tmpe99a1"
$ws5.Cells.Item(2, 12).Value = 84

# Copy the existing header/data formatting (bold + border on row 1,
# plain on row 2) across the newly added columns so the new cells match
# the look of the pre-existing B1:E2 block instead of picking up
# Excel's blank default format.
$ws5.Range("B1").Copy() | Out-Null
$ws5.Range("F1:L1").PasteSpecial(-4122) | Out-Null

$ws5.Range("B2").Copy() | Out-Null
$ws5.Range("F2:L2").PasteSpecial(-4122) | Out-Null

$ws5.Range("A1").Select() | Out-Null
